# Add data for 2022-12-04
# 1) Rename the sheet and update the header label from "Through November 25" to
#    "Through November 26" (and the matching sheet/tab name).
# 2) Update/insert several cell values across the data grid reflecting the
#    refreshed carjacking counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab) and the column B header text.
$ws.Name = "Through 2022-11-26"
$ws.Range("B1").Value = "November 2022 (through November 26)"

# Apply the cell value updates / insertions.
$updates = @{
    "CA2"  = 1
    "B5"   = 9
    "AI5"  = 6
    "AT5"  = 6
    "M6"   = 5
    "M7"   = 6
    "BE7"  = 4
    "AT10" = 3
    "B14"  = 1
    "B17"  = 3
    "BE17" = 4
    "AT19" = 1
    "M20"  = 4
    "B21"  = 1
    "X25"  = 12
    "BP25" = 4
    "X27"  = 3
    "BP32" = 1
    "CA33" = 1
    "BP41" = 4
    "CA41" = 3
    "M46"  = 2
    "AT48" = 1
    "AT50" = 2
    "AT62" = 1
    "B89"  = 1
    "BE89" = 2
    "AT90" = 1
    "M98"  = 2
    "BP99" = 1
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
